$d = $word.ActiveDocument

# The block to remove is the footer of the page: an empty paragraph, the
# "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph, and the
# "(c) 2020 ... Creative Commons Attribution" paragraph. It sits right after
# the "LOB1038: Física Experimental I (Requisito fraco)" paragraph and right
# before the trailing empty paragraphs at the end of the document.

$reqPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*LOB1038: Física Experimental I (Requisito fraco)*") {
        $reqPara = $p
    }
}
$startPara = $reqPara.Next()

$endPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Powered by Jekyll and Github pages*") {
        $endPara = $p
    }
}

$rng = $d.Range($startPara.Range.Start, $endPara.Range.End)
$rng.Delete()
